$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Twist" header column
$ws.Range("E1").Value = "Twist"

# Values for E2:E65, in row order (matches codon table rows 2-65)
$twist = @("F","T","T","F","F","T","F","F","T","T","F","F","F","F","T","F","F","T","T","F","F","F","F","T","F","F","F","F","F","F","T","F","F","F","T","T","F","F","F","T","F","F","F","T","F","F","T","F","F","T","F","F","F","F","F","F","F","T","T","F","F","T","F","F")

for ($i = 0; $i -lt $twist.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $twist[$i]
}
